# Add 2013-2014 ("_H" cycle) instruments and deployments to the NHANES
# metadata workbook, mirroring the existing 2015-2016 ("_I") / 2017-2018
# ("_J") entries.
#
# Order of operations matters for shared-string index assignment: the
# "Instruments" sheet is edited first (its new strings land at indices
# 253-264), then the "Deployments" sheet (its new strings land at
# 265-268) - this reproduces the exact shared-string table the author's
# commit produced.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Instruments sheet: insert 4 new rows, one after each existing
#    "...I-QUESTIONNAIRE" row within its group (DMQ_Family, INQ_Family,
#    DMQ, DSQ), carrying a new "_H" (2013-2014) instrument definition.
# ---------------------------------------------------------------------
$wsInstruments = $wb.Worksheets.Item("Instruments")

$wsInstruments.Range("5:5").Insert()
$wsInstruments.Cells.Item(5,1).Value2 = "nhanes-kb:INS-DMQ_Family_H-QUESTIONNAIRE"
$wsInstruments.Cells.Item(5,2).Value2 = "nhanes:00063"
$wsInstruments.Cells.Item(5,3).Value2 = "Family Questionnaire - Demographic Background (2013-2014)"
$wsInstruments.Cells.Item(5,6).Value2 = "https://wwwn.cdc.gov/nchs/data/nhanes/2013-2014/questionnaires/DMQ_Family_H.pdf"

$wsInstruments.Range("10:10").Insert()
$wsInstruments.Cells.Item(10,1).Value2 = "nhanes-kb:INS-INQ_Family_H-QUESTIONNAIRE"
$wsInstruments.Cells.Item(10,2).Value2 = "nhanes:00063"
$wsInstruments.Cells.Item(10,3).Value2 = "Family Questionnaire - Income (2013-2014)"
$wsInstruments.Cells.Item(10,6).Value2 = "https://wwwn.cdc.gov/nchs/data/nhanes/2013-2014/questionnaires/INQ_Family_H.pdf"
$wsInstruments.Rows.Item(10).RowHeight = 15

$wsInstruments.Range("19:19").Insert()
$wsInstruments.Cells.Item(19,1).Value2 = "nhanes-kb:INS-DMQ_H-QUESTIONNAIRE"
$wsInstruments.Cells.Item(19,2).Value2 = "nhanes:00064"
$wsInstruments.Cells.Item(19,3).Value2 = "Sample Person Questionnaire - Demographic (2013-2014)"
$wsInstruments.Cells.Item(19,6).Value2 = "https://wwwn.cdc.gov/nchs/data/nhanes/2013-2014/questionnaires/DMQ_H.pdf"
$wsInstruments.Rows.Item(19).RowHeight = 15

$wsInstruments.Range("25:25").Insert()
$wsInstruments.Cells.Item(25,1).Value2 = "nhanes-kb:INS-DSQ_H-QUESTIONNAIRE"
$wsInstruments.Cells.Item(25,2).Value2 = "nhanes:00064"
$wsInstruments.Cells.Item(25,3).Value2 = "Sample Person Questionnaire - Dietary Supplements and Prescription Medication (2013-2014)"
$wsInstruments.Cells.Item(25,6).Value2 = "https://wwwn.cdc.gov/nchs/data/nhanes/2013-2014/questionnaires/DSQ_H.pdf"
$wsInstruments.Rows.Item(25).RowHeight = 15

$wsInstruments.Activate()
$wsInstruments.Range("F25").Select()

# ---------------------------------------------------------------------
# 2. Deployments sheet: insert 4 new rows, one per new instrument above,
#    pairing the DPL- deployment id with the INS- instrument id.
# ---------------------------------------------------------------------
$wsDeployments = $wb.Worksheets.Item("Deployments")

$wsDeployments.Range("4:4").Insert()
$wsDeployments.Cells.Item(4,1).Value2 = "nhanes-kb:DPL-DMQ_Family_H-QUESTIONNAIRE"
$wsDeployments.Cells.Item(4,2).Value2 = "vstoi:Deployment"
$wsDeployments.Cells.Item(4,3).Value2 = "nhanes-kb:PLT-GENERIC-HUMAN"
$wsDeployments.Cells.Item(4,4).Value2 = "nhanes-kb:INS-DMQ_Family_H-QUESTIONNAIRE"
$wsDeployments.Cells.Item(4,5).Value2 = "nhanes-kb:DET-GENERIC-DETECTOR"
$wsDeployments.Cells.Item(4,6).Value2 = "2015-11-29T11:00:00.999Z"

$wsDeployments.Range("7:7").Insert()
$wsDeployments.Cells.Item(7,1).Value2 = "nhanes-kb:DPL-INQ_Family_H-QUESTIONNAIRE"
$wsDeployments.Cells.Item(7,2).Value2 = "vstoi:Deployment"
$wsDeployments.Cells.Item(7,3).Value2 = "nhanes-kb:PLT-GENERIC-HUMAN"
$wsDeployments.Cells.Item(7,4).Value2 = "nhanes-kb:INS-INQ_Family_H-QUESTIONNAIRE"
$wsDeployments.Cells.Item(7,5).Value2 = "nhanes-kb:DET-GENERIC-DETECTOR"
$wsDeployments.Cells.Item(7,6).Value2 = "2015-11-29T11:00:00.999Z"

$wsDeployments.Range("10:10").Insert()
$wsDeployments.Cells.Item(10,1).Value2 = "nhanes-kb:DPL-DMQ_H-QUESTIONNAIRE"
$wsDeployments.Cells.Item(10,2).Value2 = "vstoi:Deployment"
$wsDeployments.Cells.Item(10,3).Value2 = "nhanes-kb:PLT-GENERIC-HUMAN"
$wsDeployments.Cells.Item(10,4).Value2 = "nhanes-kb:INS-DMQ_H-QUESTIONNAIRE"
$wsDeployments.Cells.Item(10,5).Value2 = "nhanes-kb:DET-GENERIC-DETECTOR"
$wsDeployments.Cells.Item(10,6).Value2 = "2015-11-29T11:00:00.999Z"

$wsDeployments.Range("13:13").Insert()
$wsDeployments.Cells.Item(13,1).Value2 = "nhanes-kb:DPL-DSQ_H-QUESTIONNAIRE"
$wsDeployments.Cells.Item(13,2).Value2 = "vstoi:Deployment"
$wsDeployments.Cells.Item(13,3).Value2 = "nhanes-kb:PLT-GENERIC-HUMAN"
$wsDeployments.Cells.Item(13,4).Value2 = "nhanes-kb:INS-DSQ_H-QUESTIONNAIRE"
$wsDeployments.Cells.Item(13,5).Value2 = "nhanes-kb:DET-GENERIC-DETECTOR"
$wsDeployments.Cells.Item(13,6).Value2 = "2015-11-29T11:00:00.999Z"

$wsDeployments.Activate()
$wsDeployments.Range("D13").Select()
